# Add save history function
# This adds a "time" intent row and reworks the "name" intent row so the
# "my name" phrasing (and "my name is") lives together with the other
# name-recall phrases. A brand-new "change" row is appended after it to
# keep the three rows in this block intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the existing "change" row (row 12) onto the
# new row 13 so the new "time" row matches the style (s="1") of its peers.
$ws.Range("A12:B12").Copy()
$ws.Range("A13:B13").PasteSpecial(-4122)

# New row 13: the "time" intent and its trigger phrases.
$ws.Range("A13").Value = "time"
$ws.Range("B13").Value = "time, current time, what time is it, can you tell me the time, do you know the current time"

# Row 11 ("name") keeps its intent label, but the list of trigger phrases
# now also covers "my name" / "my name is" (previously these lived apart).
$ws.Range("B11").Value = "my name, do you remember my name, do you know my name, who am i, tell me who am i, call my name, what's my name, call me, my name is"

# Restore the selection to match the edited area.
$ws.Range("B11").Select()
